# A new price-record row was inserted into the daily price log as row 115,
# pushing the previously existing rows 115-178 down to 116-179 (the sheet
# grows from A1:T178 to A1:T179).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 115; existing rows 115..178 shift to 116..179.
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new record.
$ws.Cells.Item(115, 1).Value2  = 5
$ws.Cells.Item(115, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(115, 3).Value2  = "Maule"
$ws.Cells.Item(115, 4).Value2  = 44529
$ws.Cells.Item(115, 5).Value2  = 7
$ws.Cells.Item(115, 6).Value2  = "Fruta"
$ws.Cells.Item(115, 7).Value2  = 100108
$ws.Cells.Item(115, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(115, 9).Value2  = 100108005
$ws.Cells.Item(115, 10).Value2 = "Piña"
$ws.Cells.Item(115, 11).Value2 = "Caramelo"
$ws.Cells.Item(115, 12).Value2 = "Tercera"
$ws.Cells.Item(115, 13).Value2 = 200
$ws.Cells.Item(115, 14).Value2 = 18500
$ws.Cells.Item(115, 15).Value2 = 18500
$ws.Cells.Item(115, 16).Value2 = 18500
$ws.Cells.Item(115, 17).Value2 = "$/caja 16 unidades"
$ws.Cells.Item(115, 18).Value2 = "Ecuador"
$ws.Cells.Item(115, 19).Value2 = 1156
$ws.Cells.Item(115, 20).Value2 = 16
